# Updated cryptos list data (Price and Volume(1h) columns), generated from commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '65.329.39'
Set-TextCell 'E2' '  +2.49%  '
Set-TextCell 'D3' '2.655.60'
Set-TextCell 'E3' '  +1.48%  '
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '606.15'
Set-TextCell 'E5' '  +1.85%  '
Set-TextCell 'D6' '157.77'
Set-TextCell 'E6' '  +4.68%  '
Set-TextCell 'E7' '  -0.07%  '
Set-TextCell 'E9' '  +8.88%  '
Set-TextCell 'D10' '0.412'
Set-TextCell 'E10' '  +4.89%  '
Set-TextCell 'D11' '5.83'
Set-TextCell 'E11' '  +0.58%  '
Set-TextCell 'E12' '  +1.76%  '
Set-TextCell 'D13' '29.58'
Set-TextCell 'E13' '  +5.89%  '
Set-TextCell 'E14' '  +16.23%  '
Set-TextCell 'D15' '3.133.80'
Set-TextCell 'E15' '  +1.53%  '
Set-TextCell 'D16' '65.132.38'
Set-TextCell 'E16' '  +2.46%  '
Set-TextCell 'D17' '2.655.74'
Set-TextCell 'E17' '  +2.00%  '
Set-TextCell 'D18' '12.79'
Set-TextCell 'E18' '  +4.34%  '
Set-TextCell 'E19' '  +2.68%  '
Set-TextCell 'D20' '360.32'
Set-TextCell 'E20' '  +3.49%  '
Set-TextCell 'D21' '7.37'
Set-TextCell 'E21' '  +5.46%  '
Set-TextCell 'E22' '  -0.02%  '
Set-TextCell 'D23' '69.31'
Set-TextCell 'E23' '  +3.00%  '
Set-TextCell 'E24' '  +1.95%  '
Set-TextCell 'D25' '9.60'
Set-TextCell 'E25' '  +3.71%  '
Set-TextCell 'E26' '  +16.36%  '
Set-TextCell 'E27' '  -1.50%  '
Set-TextCell 'D28' '8.27'
Set-TextCell 'E28' '  -2.63%  '
Set-TextCell 'D29' '0.166'
Set-TextCell 'E29' '  +1.69%  '
Set-TextCell 'B30' 'Bittensor'
Set-TextCell 'C30' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D30' '554.41'
Set-TextCell 'E30' '  +1.00%  '
Set-TextCell 'B31' 'PancakeSwap'
Set-TextCell 'C31' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D31' '2.23'
Set-TextCell 'E31' '  +8.53%  '
Set-TextCell 'E32' '  +0.11%  '
Set-TextCell 'E33' '  +2.44%  '
Set-TextCell 'D34' '5.63'
Set-TextCell 'E34' '  +1.72%  '
Set-TextCell 'E35' '  +4.51%  '
Set-TextCell 'E36' '  +3.95%  '
Set-TextCell 'D37' '20.62'
Set-TextCell 'E37' '  +4.72%  '
Set-TextCell 'D38' '163.61'
Set-TextCell 'E38' '  -0.49%  '
Set-TextCell 'E39' '  +1.58%  '
Set-TextCell 'E40' '  -0.01%  '
Set-TextCell 'D42' '42.51'
Set-TextCell 'E42' '  +6.88%  '
Set-TextCell 'D43' '167.82'
Set-TextCell 'E43' '  +0.54%  '
Set-TextCell 'E44' '  +2.71%  '
Set-TextCell 'D45' '0.0621'
Set-TextCell 'E45' '  +6.40%  '
Set-TextCell 'D46' '2.33'
Set-TextCell 'E46' '  +8.37%  '
Set-TextCell 'D47' '23.09'
Set-TextCell 'E47' '  -2.11%  '
Set-TextCell 'D48' '0.658'
Set-TextCell 'E48' '  +3.74%  '
Set-TextCell 'E49' '  +5.17%  '
Set-TextCell 'E50' '  +2.23%  '
Set-TextCell 'D51' '19.76'
Set-TextCell 'E51' '  +2.20%  '
